$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 13:22"

# Update "Casos activos" (C) and "Muertes" (E) columns for the affected provinces

# Madrid (row 4)
$ws.Range("C4").Value = 37154
$ws.Range("E4").Value = 8222

# Cataluña (row 5)
$ws.Range("C5").Value = 19640
$ws.Range("E5").Value = 5061

# Castilla y Leon (row 6)
$ws.Range("C6").Value = 6686
$ws.Range("E6").Value = 1770

# Castilla-La Mancha (row 7)
$ws.Range("C7").Value = 5615
$ws.Range("E7").Value = 2498

# Pais Vasco (row 8)
$ws.Range("C8").Value = 11380
$ws.Range("E8").Value = 1312

# Andalucia (row 9)
$ws.Range("C9").Value = 6334
$ws.Range("E9").Value = 1238

# Galicia (row 10)
$ws.Range("C10").Value = 5816
$ws.Range("E10").Value = 555

# Aragon (row 14)
$ws.Range("C14").Value = 2382
$ws.Range("E14").Value = 749

# Navarra (row 15)
$ws.Range("C15").Value = 2185
$ws.Range("E15").Value = 548

# La Rioja (row 17)
$ws.Range("C17").Value = 2220
$ws.Range("E17").Value = 333

# Extremadura (row 23)
$ws.Range("C23").Value = 1960
$ws.Range("E23").Value = 451

# Asturias (row 30)
$ws.Range("C30").Value = 860

# Gran Canaria (row 32)
$ws.Range("C32").Value = 1151
$ws.Range("E32").Value = 136

# Cantabria (row 33)
$ws.Range("C33").Value = 1509
$ws.Range("E33").Value = 192

# Murcia (row 38)
$ws.Range("C38").Value = 1229
$ws.Range("E38").Value = 132
